$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" to H1, copying the formatting from G1 (the other header cells)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Save"

# Fill in the Save column values for rows 2-7
$saveValues = @(1, 1, 0, 0, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
